$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "otro usuario"
$ws.Range("D5").Select()
